$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.466.72"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "2.107.51"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.81%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.11"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4556"
$ws.Range("E8").Value = "  +5.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.46"
$ws.Range("E9").Value = "  +16.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08942"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.30"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "2.119.10"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.868"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.049"
$ws.Range("E15").Value = "  +5.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.51"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001154"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06661"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.22"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.356"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "30.525.82"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("D26").Value = "2.374.02"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.68"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.545"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.99"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.222"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.654"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.361"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.949"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.49"
$ws.Range("E36").Value = "  +8.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02582"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.677"
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06840"
$ws.Range("E39").Value = "  +3.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2305"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6879"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.250"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.007"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.331"
$ws.Range("E45").Value = "  +6.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.12"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6373"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.673"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.249"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3409"
$ws.Range("E50").Value = "  +25.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.34"
$ws.Range("E51").Value = "  +2.62%  "
